# Fruta / hortaliza, semanal
# A new weekly price-report row is inserted at row 203 (Florida King / Primera,
# Provincia de Limari), pushing the existing rows 203-211 down to 204-212.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 203, shifting rows 203:211
# down to 204:212 (matches the dimension growing from A1:T211 to A1:T212).
$ws.Rows(203).Insert()

# Populate the newly inserted row 203 with the new market record.
$ws.Cells.Item(203, 1).Value = 4
$ws.Cells.Item(203, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(203, 3).Value = 'Los Lagos'
$ws.Cells.Item(203, 4).Value = 44890
$ws.Cells.Item(203, 5).Value = 10
$ws.Cells.Item(203, 6).Value = 'Fruta'
$ws.Cells.Item(203, 7).Value = 100103
$ws.Cells.Item(203, 8).Value = 'Frutos de hueso (carozo)'
$ws.Cells.Item(203, 9).Value = 100103004
$ws.Cells.Item(203, 10).Value = 'Durazno'
$ws.Cells.Item(203, 11).Value = 'Florida King'
$ws.Cells.Item(203, 12).Value = 'Primera'
$ws.Cells.Item(203, 13).Value = 400
$ws.Cells.Item(203, 14).Value = 15000
$ws.Cells.Item(203, 15).Value = 16000
$ws.Cells.Item(203, 16).Value = 15500
$ws.Cells.Item(203, 17).Value = '$/bandeja 10 kilos empedrada'
$ws.Cells.Item(203, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(203, 19).Value = 1550
$ws.Cells.Item(203, 20).Value = 10
